# 3.1.8 chapter cleanup: remove the stray "_GoBack" bookmark that Word
# leaves behind after an editing session. It wraps no visible text (it
# sits right after "...holder som regel her. " at the end of the
# paragraph), so it is invisible to a normal text Find/Replace and has
# to be removed as an actual bookmark object.
$d = $word.ActiveDocument

# "_GoBack" is a hidden bookmark (leading underscore) so it does not show
# up in $d.Bookmarks, but it can still be addressed directly by name.
$goBack = $d.Bookmarks("_GoBack")
if ($goBack -ne $null) {
    $goBack.Delete()
}
